$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 146.925420661
$ws.Range("C2").Value = 0.000000289581149031

$ws.Range("B3").Value = 146925.420661
$ws.Range("C3").Value = 0.00289581149031

$ws.Range("B4").Value = 293850.841322
$ws.Range("C4").Value = 0.00868743447093

$ws.Range("B5").Value = 5877.016826440001
$ws.Range("C5").Value = 0.0001737486894186
